# Add team record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the existing header row (bold, bordered,
# centered). Copy format from the last existing header cell (AC1) first,
# then overwrite the values with the new header labels.
$ws.Range("AC1").Copy()
$headerRange = $ws.Range("AD1:AF1")
$headerRange.PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-50) gets the same team record: 71 wins, 91 losses, 0 ties.
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
